$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FLAM")
$ws.Name = "FLAME"
$ws.Range("A3").Value = "FLAME"
$ws.Activate()
$ws.Range("B11").Select()
